$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resource")

# Shift header labels left starting at column F (6) and rename a couple along the way,
# then clear the now-unused trailing column (N).
$ws.Range("F1").Value = "units"
$ws.Range("G1").Value = "cost_per_unit"
$ws.Range("H1").Value = "parts"
$ws.Range("I1").Value = "capacity"
$ws.Range("J1").Value = "holding_capacity"
$ws.Range("K1").Value = "skills"
$ws.Range("L1").Value = "aggregates"
$ws.Range("M1").Value = "kwargs"

# N1 ("kwargs") no longer exists after the shift - clear it entirely so the
# sheet's used range/dimension shrinks back to B1:M1.
$ws.Range("N1").Clear()
